# Rotate the observation records in rows 3, 4, 6 and 7:
#   new row 3 <- old row 7
#   new row 4 <- old row 3
#   new row 6 <- old row 4
#   new row 7 <- old row 6
# (row 5 is untouched)
#
# Columns A,B,E,F,G,H,Q,R carry the record data that differs between the
# rows; the remaining columns (C,D,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY)
# happen to already hold identical values across these four rows, so they
# do not need to be touched. A handful of columns (J,L,M,AF) only ever
# contain an empty placeholder in some of the rows; those are cleared /
# recreated to match the row the data moved from.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","E","F","G","H","Q","R")

# Snapshot the current ("before") values of the four affected rows before
# overwriting anything.
$orig = @{}
foreach ($r in 3,4,6,7) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $row
}

# Apply the rotation.
$mapping = @{ 3 = 7; 4 = 3; 6 = 4; 7 = 6 }
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $orig[$srcRow][$c]
    }
}

# The placeholder/optional cells (empty "blank" markers) travel with the
# row too: row 3 gains an empty L, row 4 loses its empty L, row 6 gains an
# empty J plus an empty AF (and loses its "färska spår" M), row 7 loses its
# empty J and AF but gains the "färska spår" M text.
# A lone "'" produces Excel's empty-text-cell (like the other blank marker
# cells, e.g. I3/J3/K3) rather than clearing the cell outright; reset the
# style afterwards so the implicit "quote prefix" formatting Excel applies
# doesn't stick around (the source blank cells carry no special style).
function Set-EmptyTextCell($addr) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

Set-EmptyTextCell "L3"
$ws.Range("L4").ClearContents()

Set-EmptyTextCell "J6"
Set-EmptyTextCell "AF6"
$ws.Range("M6").ClearContents()

$ws.Range("J7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("M7").Value = "färska spår"
